$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11 "Marking" changes
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total" changes
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "64 / 112"
